# Natmi following Dr Hou advice
#
# The Dlk2 -> Notch1 ligand/receptor edge table previously mixed up the
# "ECs" and "sCs" sending-cluster labels, and only reported one row per
# target cluster. This refreshes rows 2-4 (ECs sending cluster x ECs/FAPs/sCs
# targets) with corrected statistics and adds rows 5-7 (sCs sending cluster x
# ECs/FAPs/sCs targets) that were missing before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Notch1 (via Dlk2) -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dlk2"
$ws.Range("C2").Value = "Notch1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1179396666666667
$ws.Range("H2").Value = 0.353819
$ws.Range("I2").Value = 0.1771592059007171
$ws.Range("J2").Value = 0.1771592059007171
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 55.908252
$ws.Range("N2").Value = 167.724756
$ws.Range("O2").Value = 0.6412441619121594
$ws.Range("P2").Value = 0.6412441619121594
$ws.Range("Q2").Value = 6.593800604796001
$ws.Range("R2").Value = 59.34420544316401
$ws.Range("S2").Value = 0.113602306512829
$ws.Range("T2").Value = 0.113602306512829

# Row 3: ECs -> Notch1 (via Dlk2) -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dlk2"
$ws.Range("C3").Value = "Notch1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1179396666666667
$ws.Range("H3").Value = 0.353819
$ws.Range("I3").Value = 0.1771592059007171
$ws.Range("J3").Value = 0.1771592059007171
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.370676
$ws.Range("N3").Value = 16.112028
$ws.Range("O3").Value = 0.06159939735768789
$ws.Range("P3").Value = 0.06159939735768789
$ws.Range("Q3").Value = 0.6334157372146667
$ws.Range("R3").Value = 5.700741634932001
$ws.Range("S3").Value = 0.01091290031985072
$ws.Range("T3").Value = 0.01091290031985072

# Row 4: ECs -> Notch1 (via Dlk2) -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dlk2"
$ws.Range("C4").Value = "Notch1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1179396666666667
$ws.Range("H4").Value = 0.353819
$ws.Range("I4").Value = 0.1771592059007171
$ws.Range("J4").Value = 0.1771592059007171
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 25.90822366666667
$ws.Range("N4").Value = 77.724671
$ws.Range("O4").Value = 0.2971564407301527
$ws.Range("P4").Value = 0.2971564407301527
$ws.Range("Q4").Value = 3.055607263172111
$ws.Range("R4").Value = 27.500465368549
$ws.Range("S4").Value = 0.05264399906803735
$ws.Range("T4").Value = 0.05264399906803735

# Row 5: sCs -> Notch1 (via Dlk2) -> ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Dlk2"
$ws.Range("C5").Value = "Notch1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5477873333333333
$ws.Range("H5").Value = 1.643362
$ws.Range("I5").Value = 0.8228407940992829
$ws.Range("J5").Value = 0.822840794099283
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 55.908252
$ws.Range("N5").Value = 167.724756
$ws.Range("O5").Value = 0.6412441619121594
$ws.Range("P5").Value = 0.6412441619121594
$ws.Range("Q5").Value = 30.625832274408
$ws.Range("R5").Value = 275.632490469672
$ws.Range("S5").Value = 0.5276418553993304
$ws.Range("T5").Value = 0.5276418553993304

# Row 6: sCs -> Notch1 (via Dlk2) -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Dlk2"
$ws.Range("C6").Value = "Notch1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5477873333333333
$ws.Range("H6").Value = 1.643362
$ws.Range("I6").Value = 0.8228407940992829
$ws.Range("J6").Value = 0.822840794099283
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.370676
$ws.Range("N6").Value = 16.112028
$ws.Range("O6").Value = 0.06159939735768789
$ws.Range("P6").Value = 0.06159939735768789
$ws.Range("Q6").Value = 2.941988284237333
$ws.Range("R6").Value = 26.477894558136
$ws.Range("S6").Value = 0.05068649703783717
$ws.Range("T6").Value = 0.05068649703783718

# Row 7: sCs -> Notch1 (via Dlk2) -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Dlk2"
$ws.Range("C7").Value = "Notch1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5477873333333333
$ws.Range("H7").Value = 1.643362
$ws.Range("I7").Value = 0.8228407940992829
$ws.Range("J7").Value = 0.822840794099283
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.90822366666667
$ws.Range("N7").Value = 77.724671
$ws.Range("O7").Value = 0.2971564407301527
$ws.Range("P7").Value = 0.2971564407301527
$ws.Range("Q7").Value = 14.19219675376689
$ws.Range("R7").Value = 127.729770783902
$ws.Range("S7").Value = 0.2445124416621153
$ws.Range("T7").Value = 0.2445124416621154
